# close #77 eliminate the problem of wrong-encoding of dig form
# Append the three new error-code rows to the "Error" sheet / table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 2000
$ws.Range("B16").Value = "需要等级到达{0}级"

$ws.Range("A17").Value = 2001
$ws.Range("B17").Value = "体力不足"

$ws.Range("A18").Value = 3000
$ws.Range("B18").Value = "钻石不足"

# Grow the worksheet table ("表1") so the new rows belong to it, like
# Excel does automatically when you type into the row right below a table.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B18"))

# Match the scrolled viewport / selection from the authored workbook.
$ws.Range("B18").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
